# Generate Report for Handback
# This script mirrors a "handback" report regeneration that:
#  - Marks the de-de language as handed back (new handback datetime)
#  - Refreshes the zh-cn handback datetime placeholder text as well
#  - Populates the (previously empty) "Latest Target File" / "Latest Handback File"
#    columns (F/G) for both the zh-cn and de-de detail sheets, with hyperlinks
#  - Updates the Status text on the Overview sheet

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text change (affects Overview!B2,C2,B3,C3 via the shared string)
# ---------------------------------------------------------------------------
$wsOverview.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")

# ---------------------------------------------------------------------------
# 2. Handback datetime placeholders
#    zh-cn keeps the same "not handed back yet" cells but the placeholder
#    text itself is refreshed; de-de actually receives a real handback time.
# ---------------------------------------------------------------------------
$wsZh.Cells.Replace("0001-01-01 00:00:00", "2016-03-11 08:41:28")
$wsDe.Cells.Replace("0001-01-01 00:00:00", "2016-03-11 08:41:33")

# ---------------------------------------------------------------------------
# Helper URLs reused for the new hyperlinks (same targets as the existing
# "Latest Handoff File" (A) / "Latest Target File" (D) hyperlinks).
# ---------------------------------------------------------------------------
$mdUrlUuid1 = "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/5efaea19-8b46-491c-b5f0-3fc48d58ff97.md"
$mdUrlUuid2 = "https://github.com/OpenLocalizationTest/oltest/blob/da7a1e2757d4465740352157035a46a90c088d33/e2e/b351d048-6e0c-47da-b175-fa0e1a3d2857.md"

$zhXlfUrlUuid1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5faede3304f846602d1cca11ed3f74baec1c148/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5efaea19-8b46-491c-b5f0-3fc48d58ff97.471316cc5dac540c8bf5e51d30016e21aaaa470d.zh-cn.xlf"
$zhXlfUrlUuid2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5faede3304f846602d1cca11ed3f74baec1c148/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b351d048-6e0c-47da-b175-fa0e1a3d2857.44d9580c31a1f8dbc94519769a85fd88329c986e.zh-cn.xlf"

$deXlfUrlUuid1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69340c41d8f5df6026ddcbbe2603c08d4516ede0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5efaea19-8b46-491c-b5f0-3fc48d58ff97.471316cc5dac540c8bf5e51d30016e21aaaa470d.de-de.xlf"
$deXlfUrlUuid2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69340c41d8f5df6026ddcbbe2603c08d4516ede0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b351d048-6e0c-47da-b175-fa0e1a3d2857.44d9580c31a1f8dbc94519769a85fd88329c986e.de-de.xlf"

# Display text for the new cells.
$mdDisplay1 = "5efaea19-8b46-491c-b5f0-3fc48d58ff97.md"
$mdDisplay2 = "b351d048-6e0c-47da-b175-fa0e1a3d2857.md"
$zhXlfDisplay1 = "5efaea19-8b46-491c-b5f0-3fc48d58ff97.471316cc5dac540c8bf5e51d30016e21aaaa470d.zh-cn.xlf"
$zhXlfDisplay2 = "b351d048-6e0c-47da-b175-fa0e1a3d2857.44d9580c31a1f8dbc94519769a85fd88329c986e.zh-cn.xlf"
$deXlfDisplay1 = "5efaea19-8b46-491c-b5f0-3fc48d58ff97.471316cc5dac540c8bf5e51d30016e21aaaa470d.de-de.xlf"
$deXlfDisplay2 = "b351d048-6e0c-47da-b175-fa0e1a3d2857.44d9580c31a1f8dbc94519769a85fd88329c986e.de-de.xlf"

function Add-HandbackLink($ws, $cellRef, $text, $url) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $text
    $ws.Hyperlinks.Add($cell, $url, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
    # Match the look of the existing "hyperlink style" cells (A/B/D columns):
    # underlined, blue font - same as the workbook's built-in Hyperlink style.
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# 3. Populate the new F (Latest Target File) / G (Latest Handback File)
#    columns for the zh-cn sheet, rows 2 and 3.
# ---------------------------------------------------------------------------
Add-HandbackLink $wsZh "F2" $mdDisplay1 $mdUrlUuid1
Add-HandbackLink $wsZh "G2" $zhXlfDisplay1 $zhXlfUrlUuid1
Add-HandbackLink $wsZh "F3" $mdDisplay2 $mdUrlUuid2
Add-HandbackLink $wsZh "G3" $zhXlfDisplay2 $zhXlfUrlUuid2

# ---------------------------------------------------------------------------
# 4. Populate the new F / G columns for the de-de sheet, rows 2 and 3.
# ---------------------------------------------------------------------------
Add-HandbackLink $wsDe "F2" $mdDisplay1 $mdUrlUuid1
Add-HandbackLink $wsDe "G2" $deXlfDisplay1 $deXlfUrlUuid1
Add-HandbackLink $wsDe "F3" $mdDisplay2 $mdUrlUuid2
Add-HandbackLink $wsDe "G3" $deXlfDisplay2 $deXlfUrlUuid2

Write-Host "Handback report generated."
